$wb = $excel.ActiveWorkbook

# --- Fix header typos on the AddCustomerTest sheet ---
$ws1 = $wb.Worksheets.Item("AddCustomerTest")
$ws1.Range("A1").Value = "firstName"
$ws1.Range("E1").Value = "runmode"

# --- Set the run mode for OpenAccountTest to Y on the test_suite sheet ---
$ws3 = $wb.Worksheets.Item("test_suite")
$ws3.Range("B4").Value = "Y"

# --- Reflect the user's last selections / active sheet ---
$ws3.Activate()
$ws3.Range("B4").Select()

$ws1.Activate()
$ws1.Range("E3").Select()
